# Refresh the Victorian "key outbreaks" cluster table with the newer
# snapshot: add newly listed clusters, drop resolved ones, update the
# active-case counts, and flip each row's Exist marker (old/new).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('126 Racecourse Road Public Housing Tower Flemington', 10, 'new'),
    @('3175 The Bays Aged Care Facility Hastings', 14, 'old'),
    @('3175 The Bays Aged Care Facility Hastings', 17, 'new'),
    @('3563 Embracia Aged Care Reservoir', 22, 'old'),
    @('3612 BlueCross Glengowrie Outbreak', 10, 'new'),
    @('3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak', 20, 'new'),
    @('4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak', 14, 'new'),
    @('Apartment Complex Fawkner', 10, 'old'),
    @('Australian Lamb Colac East', 13, 'old'),
    @('Australian Lamb Colac East', 14, 'new'),
    @('Bread Solutions Braeside', 13, 'old'),
    @('Bread Solutions Braeside Outbreak', 20, 'new'),
    @('CS Square Caroline Springs', 13, 'old'),
    @('CS Square Caroline Springs', 15, 'new'),
    @('Carton Finishing Pty. Ltd. Campbellfield', 10, 'new'),
    @('Carton Finishing Pty. Ltd. Campbellfield', 12, 'old'),
    @('Cedar Meats Australia Brooklyn', 10, 'old'),
    @('Cedar Meats Australia Brooklyn Outbreak', 11, 'new'),
    @('Community Kids Bayswater Early Education Centre Bayswater North', 10, 'new'),
    @('Community Kids Bayswater Early Education Centre Bayswater North', 18, 'old'),
    @('Costco Wholesale Epping', 13, 'old'),
    @('Embracia Aged Care Reservoir Outbreak', 23, 'new'),
    @('Ermha365 Residential Disability Care Services Doveton', 10, 'old'),
    @('FedEx Station Melbourne Airport', 11, 'new'),
    @('FedEx Station Melbourne Airport', 14, 'old'),
    @('Green Leaves Early Learning Centre Highlands Craigieburn', 10, 'new'),
    @('Green Leaves Early Learning Centre Highlands Craigieburn', 14, 'old'),
    @('Guardian Childcare Caulfield', 14, 'old'),
    @('Guardian Childcare Caulfield Outbreak', 15, 'new'),
    @('Kool Kidz Childcare Narre Warren', 11, 'new'),
    @('Kool Kidz Childcare Narre Warren', 16, 'old'),
    @('Melbourne Assessment Prison West Melbourne', 10, 'old'),
    @('MyCentre Childcare Broadmeadows', 10, 'old'),
    @('Nido Early School Ascot Vale', 15, 'new'),
    @('Nido Early School Ascot Vale', 29, 'old'),
    @('Nido Early School Glenroy', 23, 'new'),
    @('Nido Early School Glenroy', 24, 'old'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 48, 'new'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 52, 'old'),
    @('Northern Health The Northern Hospital Epping', 16, 'new'),
    @('Northern Health The Northern Hospital Epping', 21, 'old'),
    @('Oceania Meat Processors Laverton North Outbreak', 13, 'new'),
    @('Pick It Up Fitness Mulgrave Outbreak', 11, 'new'),
    @('Social Gathering Warrnambool 28 Sep Outbreak', 17, 'new'),
    @('Social Gathering Warrnambool 28 September', 17, 'old'),
    @('St Margaret''s Primary School OSHC Maribyrnong', 11, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne', 35, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne', 37, 'new'),
    @('Target Distribution Centre Truganina Outbreak', 21, 'new'),
    @('The Robin Hood Inn Drouin West Outbreak', 16, 'new'),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 13, 'new'),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 16, 'old'),
    @('The Royal Talbot Rehabilitation Centre', 12, 'old'),
    @('Wallaby Childcare Wollert', 16, 'old'),
    @('Werribee Mercy Hospital Emergency Department', 20, 'old'),
    @('Werribee Mercy Hospital Emergency Department', 22, 'new'),
    @('Western Health Footscray Hospital Emergency Department', 10, 'old'),
    @('Western Health Sunshine Hospital Emergency Department', 16, 'old'),
    @('Western Health Sunshine Hospital Emergency Department', 18, 'new'),
)

$rowCount = $data.Count
$lastRow = 1 + $rowCount

# Clear any previously used rows below the header before rewriting them
$usedRows = $ws.UsedRange.Rows.Count
if ($usedRows -gt 1) {
    $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($usedRows, 3)).ClearContents()
}

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
